$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USE_ACTUAL_MODEL column (F): the rows that fed model "MLB" now
# reference the optimizer configuration "PR_B_Y2/PR_B_Y3", and the rows that
# fed model "MLA" now reference "PR_B_Y3".
$ws.Range("F2:F6").Value = "PR_B_Y2/PR_B_Y3"
$ws.Range("F7:F9").Value = "PR_B_Y3"

# Move/save the selection cursor to H6, matching the saved view state.
$ws.Range("H6").Select()
